# "changed the infinity loop, added contact in schedule file"
#
# Net semantic change (per the supplied OOXML diff):
#   - A brand-new slide titled "Optional" is inserted right after the
#     existing "Recursion" slide (old/new index 36) and right before
#     "Recursion: Example" (old index 37). Everything from "Recursion:
#     Example" onward simply shifts down by one slide; none of their
#     content changes.
#   - The new slide has a single Title placeholder (no body/content
#     placeholder at all), repositioned/resized to:
#       off  x=838200  y=2514600   (EMU)
#       ext  cx=7772400 cy=1143000 (EMU)
#     and its text is "Optional".
#
# PowerPoint's object model works in points, so EMU -> pt is EMU/12700.

$p = $ppt.ActivePresentation

$recursionSlideIndex = 36

# Same master as the rest of the deck; "Title Only" layout (#6) has just
# a title placeholder defined, matching the single <p:sp> the new slide
# ends up with (no leftover/ghost content placeholder).
$titleOnlyLayout = $p.SlideMaster.CustomLayouts.Item(6)

$newSlide = $p.Slides.AddSlide($recursionSlideIndex + 1, $titleOnlyLayout)

$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Optional"

$title.Left   = 838200 / 12700
$title.Top    = 2514600 / 12700
$title.Width  = 7772400 / 12700
$title.Height = 1143000 / 12700
